# Generate Report for Handoff
# - Updates the "Priority" column (E) to "ht" for the rows whose latest
#   handback datetime matches the most-recent handoff run, on both the
#   zh-cn and de-de language sheets.
# - Bumps the "Latest Handback DateTime" / "Latest HO Xliff Generate Date"
#   timestamps to reflect the new report generation run.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

# Rows (in the zh-cn / de-de tables) whose "Latest Handback DateTime"
# equals the timestamp that is being refreshed by this report run.
$rows = @(7, 9, 10, 11, 13, 14)

foreach ($r in $rows) {
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"
}

# zh-cn sheet: "Latest Handback DateTime" (column H) moves forward.
foreach ($r in $rows) {
    $zhcn.Range("H$r").Value = "2016-08-31 06:23:59"
}

# de-de sheet: "Latest Handback DateTime" (column H) moves forward.
foreach ($r in $rows) {
    $dede.Range("H$r").Value = "2016-08-31 06:24:13"
}

# Overview sheet: "Latest HO Xliff Generate Date" (column G) mirrors the
# de-de timestamp update.
foreach ($r in $rows) {
    $overview.Range("G$r").Value = "2016-08-31 06:24:13"
}
